$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Name = "jscosc"

$ws.Cells.Item(1, 10).Value = 31.18170571327209
$ws.Cells.Item(2, 2).Value = 1864
$ws.Cells.Item(2, 4).Value = 1863
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 99.94635193133047
$ws.Cells.Item(2, 8).Value = 100
$ws.Cells.Item(2, 10).Value = 42.21394276618958
$ws.Cells.Item(3, 2).Value = 2083
$ws.Cells.Item(3, 4).Value = 2082
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 99.951992318771
$ws.Cells.Item(3, 9).Value = 0.0004798464491362764
$ws.Cells.Item(3, 10).Value = 38.75256586074829
$ws.Cells.Item(4, 2).Value = 2592
$ws.Cells.Item(4, 4).Value = 2564
$ws.Cells.Item(4, 5).Value = 27
$ws.Cells.Item(4, 6).Value = 7
$ws.Cells.Item(4, 7).Value = 99.72773239984441
$ws.Cells.Item(4, 8).Value = 98.95793130065611
$ws.Cells.Item(4, 9).Value = 0.01321928460342146
$ws.Cells.Item(4, 10).Value = 35.46834683418274
$ws.Cells.Item(5, 2).Value = 2052
$ws.Cells.Item(5, 4).Value = 2008
$ws.Cells.Item(5, 5).Value = 43
$ws.Cells.Item(5, 6).Value = 18
$ws.Cells.Item(5, 7).Value = 99.11154985192498
$ws.Cells.Item(5, 8).Value = 97.90346172598733
$ws.Cells.Item(5, 9).Value = 0.03009373458312778
$ws.Cells.Item(5, 10).Value = 43.96979928016663
$ws.Cells.Item(6, 4).Value = 1755
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 7
$ws.Cells.Item(6, 7).Value = 99.60272417707151
$ws.Cells.Item(6, 8).Value = 99.82935153583618
$ws.Cells.Item(6, 9).Value = 0.005672149744753261
$ws.Cells.Item(6, 10).Value = 33.63885712623596
$ws.Cells.Item(7, 10).Value = 37.47630786895752
$ws.Cells.Item(8, 2).Value = 2123
$ws.Cells.Item(8, 4).Value = 2122
$ws.Cells.Item(8, 5).Value = 0
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 99.95289684408856
$ws.Cells.Item(8, 8).Value = 100
$ws.Cells.Item(8, 9).Value = 0.0004708097928436911
$ws.Cells.Item(8, 10).Value = 33.0213565826416
$ws.Cells.Item(9, 10).Value = 35.85972547531128
$ws.Cells.Item(10, 2).Value = 3216
$ws.Cells.Item(10, 4).Value = 1794
$ws.Cells.Item(10, 5).Value = 1422
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(10, 7).Value = 100
$ws.Cells.Item(10, 8).Value = 55.78358208955224
$ws.Cells.Item(10, 9).Value = 0.792200557103064
$ws.Cells.Item(10, 10).Value = 43.79690098762512
$ws.Cells.Item(11, 2).Value = 1863
$ws.Cells.Item(11, 4).Value = 1861
$ws.Cells.Item(11, 6).Value = 17
$ws.Cells.Item(11, 7).Value = 99.09478168264111
$ws.Cells.Item(11, 8).Value = 99.94629430719657
$ws.Cells.Item(11, 9).Value = 0.009579563597658328
$ws.Cells.Item(11, 10).Value = 29.35353708267212
$ws.Cells.Item(12, 10).Value = 32.65231513977051
$ws.Cells.Item(13, 2).Value = 2388
$ws.Cells.Item(13, 4).Value = 2387
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 24
$ws.Cells.Item(13, 7).Value = 99.00456242223144
$ws.Cells.Item(13, 8).Value = 100
$ws.Cells.Item(13, 9).Value = 0.009950248756218905
$ws.Cells.Item(13, 10).Value = 31.61697912216187
$ws.Cells.Item(14, 10).Value = 31.17237329483032
$ws.Cells.Item(15, 2).Value = 2281
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 8).Value = 99.86842105263158
$ws.Cells.Item(15, 9).Value = 0.00131694468832309
$ws.Cells.Item(15, 10).Value = 33.32834458351135
$ws.Cells.Item(16, 2).Value = 1996
$ws.Cells.Item(16, 5).Value = 9
$ws.Cells.Item(16, 8).Value = 99.54887218045113
$ws.Cells.Item(16, 9).Value = 0.004529441368897836
$ws.Cells.Item(16, 10).Value = 36.52176094055176
$ws.Cells.Item(17, 10).Value = 36.59190249443054
$ws.Cells.Item(18, 10).Value = 35.74937438964844
$ws.Cells.Item(19, 10).Value = 36.34618973731995
$ws.Cells.Item(20, 10).Value = 28.22540903091431
$ws.Cells.Item(21, 2).Value = 2600
$ws.Cells.Item(21, 5).Value = 1
$ws.Cells.Item(21, 8).Value = 99.96152366294729
$ws.Cells.Item(21, 9).Value = 0.001153402537485583
$ws.Cells.Item(21, 10).Value = 34.35290789604187
$ws.Cells.Item(22, 2).Value = 1944
$ws.Cells.Item(22, 4).Value = 1943
$ws.Cells.Item(22, 6).Value = 19
$ws.Cells.Item(22, 7).Value = 99.03160040774719
$ws.Cells.Item(22, 9).Value = 0.00967906265919511
$ws.Cells.Item(22, 10).Value = 41.34923934936523
$ws.Cells.Item(23, 2).Value = 2040
$ws.Cells.Item(23, 4).Value = 2037
$ws.Cells.Item(23, 5).Value = 2
$ws.Cells.Item(23, 6).Value = 98
$ws.Cells.Item(23, 7).Value = 95.40983606557377
$ws.Cells.Item(23, 8).Value = 99.90191270230505
$ws.Cells.Item(23, 9).Value = 0.04681647940074907
$ws.Cells.Item(23, 10).Value = 33.44520044326782
$ws.Cells.Item(24, 2).Value = 2947
$ws.Cells.Item(24, 4).Value = 2923
$ws.Cells.Item(24, 5).Value = 23
$ws.Cells.Item(24, 6).Value = 56
$ws.Cells.Item(24, 7).Value = 98.12017455521988
$ws.Cells.Item(24, 8).Value = 99.21928038017651
$ws.Cells.Item(24, 9).Value = 0.02651006711409396
$ws.Cells.Item(24, 10).Value = 41.00025534629822
$ws.Cells.Item(25, 2).Value = 2648
$ws.Cells.Item(25, 4).Value = 2647
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 6).Value = 8
$ws.Cells.Item(25, 7).Value = 99.69868173258004
$ws.Cells.Item(25, 8).Value = 100
$ws.Cells.Item(25, 9).Value = 0.003012048192771084
$ws.Cells.Item(25, 10).Value = 34.70510077476501
$ws.Cells.Item(26, 2).Value = 1851
$ws.Cells.Item(26, 4).Value = 1846
$ws.Cells.Item(26, 5).Value = 4
$ws.Cells.Item(26, 6).Value = 13
$ws.Cells.Item(26, 7).Value = 99.30069930069931
$ws.Cells.Item(26, 8).Value = 99.78378378378379
$ws.Cells.Item(26, 9).Value = 0.009139784946236559
$ws.Cells.Item(26, 10).Value = 37.22248387336731
$ws.Cells.Item(27, 2).Value = 2948
$ws.Cells.Item(27, 4).Value = 2940
$ws.Cells.Item(27, 5).Value = 7
$ws.Cells.Item(27, 6).Value = 14
$ws.Cells.Item(27, 7).Value = 99.52606635071091
$ws.Cells.Item(27, 8).Value = 99.76247030878859
$ws.Cells.Item(27, 9).Value = 0.007106598984771574
$ws.Cells.Item(27, 10).Value = 37.85532569885254
$ws.Cells.Item(28, 10).Value = 34.9201328754425
$ws.Cells.Item(29, 2).Value = 2631
$ws.Cells.Item(29, 4).Value = 2627
$ws.Cells.Item(29, 5).Value = 3
$ws.Cells.Item(29, 6).Value = 22
$ws.Cells.Item(29, 7).Value = 99.16949792374481
$ws.Cells.Item(29, 8).Value = 99.88593155893535
$ws.Cells.Item(29, 9).Value = 0.009433962264150943
$ws.Cells.Item(29, 10).Value = 35.28027582168579
$ws.Cells.Item(30, 10).Value = 33.548255443573
$ws.Cells.Item(31, 4).Value = 3248
$ws.Cells.Item(31, 5).Value = 1
$ws.Cells.Item(31, 6).Value = 2
$ws.Cells.Item(31, 7).Value = 99.93846153846154
$ws.Cells.Item(31, 8).Value = 99.96922129886119
$ws.Cells.Item(31, 9).Value = 0.0009227929867733005
$ws.Cells.Item(31, 10).Value = 35.27323198318481
$ws.Cells.Item(32, 2).Value = 2273
$ws.Cells.Item(32, 4).Value = 2257
$ws.Cells.Item(32, 5).Value = 15
$ws.Cells.Item(32, 6).Value = 4
$ws.Cells.Item(32, 7).Value = 99.82308712958867
$ws.Cells.Item(32, 8).Value = 99.33978873239437
$ws.Cells.Item(32, 9).Value = 0.008399646330680813
$ws.Cells.Item(32, 10).Value = 37.08492207527161
$ws.Cells.Item(33, 2).Value = 3362
$ws.Cells.Item(33, 4).Value = 3361
$ws.Cells.Item(33, 6).Value = 1
$ws.Cells.Item(33, 7).Value = 99.97025580011898
$ws.Cells.Item(33, 9).Value = 0.0002973535533749628
$ws.Cells.Item(33, 10).Value = 36.52109241485596
$ws.Cells.Item(34, 2).Value = 2154
$ws.Cells.Item(34, 4).Value = 2153
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(34, 7).Value = 100
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 10).Value = 36.86783218383789
$ws.Cells.Item(35, 10).Value = 42.01372909545898
$ws.Cells.Item(36, 2).Value = 2426
$ws.Cells.Item(36, 4).Value = 2416
$ws.Cells.Item(36, 5).Value = 9
$ws.Cells.Item(36, 6).Value = 10
$ws.Cells.Item(36, 7).Value = 99.58779884583677
$ws.Cells.Item(36, 8).Value = 99.62886597938144
$ws.Cells.Item(36, 9).Value = 0.007828594973217964
$ws.Cells.Item(36, 10).Value = 36.01111388206482
$ws.Cells.Item(37, 2).Value = 2411
$ws.Cells.Item(37, 4).Value = 2410
$ws.Cells.Item(37, 5).Value = 0
$ws.Cells.Item(37, 6).Value = 72
$ws.Cells.Item(37, 7).Value = 97.09911361804996
$ws.Cells.Item(37, 8).Value = 100
$ws.Cells.Item(37, 10).Value = 39.30692028999329
$ws.Cells.Item(38, 2).Value = 2606
$ws.Cells.Item(38, 4).Value = 2604
$ws.Cells.Item(38, 5).Value = 1
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(38, 7).Value = 100
$ws.Cells.Item(38, 8).Value = 99.9616122840691
$ws.Cells.Item(38, 9).Value = 0.0003838771593090211
$ws.Cells.Item(38, 10).Value = 31.53883671760559
$ws.Cells.Item(39, 2).Value = 2056
$ws.Cells.Item(39, 4).Value = 2050
$ws.Cells.Item(39, 5).Value = 5
$ws.Cells.Item(39, 6).Value = 2
$ws.Cells.Item(39, 7).Value = 99.90253411306043
$ws.Cells.Item(39, 8).Value = 99.75669099756691
$ws.Cells.Item(39, 9).Value = 0.003409644422795908
$ws.Cells.Item(39, 10).Value = 37.33768248558044
$ws.Cells.Item(40, 10).Value = 37.95244574546814
$ws.Cells.Item(41, 10).Value = 35.5672972202301
$ws.Cells.Item(42, 2).Value = 1779
$ws.Cells.Item(42, 5).Value = 1
$ws.Cells.Item(42, 8).Value = 99.9437570303712
$ws.Cells.Item(42, 9).Value = 0.001685393258426966
$ws.Cells.Item(42, 10).Value = 34.55413031578064
$ws.Cells.Item(43, 2).Value = 3078
$ws.Cells.Item(43, 4).Value = 3071
$ws.Cells.Item(43, 5).Value = 6
$ws.Cells.Item(43, 6).Value = 7
$ws.Cells.Item(43, 7).Value = 99.772579597141
$ws.Cells.Item(43, 8).Value = 99.80500487487812
$ws.Cells.Item(43, 9).Value = 0.004222150048717116
$ws.Cells.Item(43, 10).Value = 41.9824230670929
$ws.Cells.Item(44, 2).Value = 2753
$ws.Cells.Item(44, 4).Value = 2752
$ws.Cells.Item(44, 5).Value = 0
$ws.Cells.Item(44, 6).Value = 0
$ws.Cells.Item(44, 7).Value = 100
$ws.Cells.Item(44, 8).Value = 100
$ws.Cells.Item(44, 9).Value = 0
$ws.Cells.Item(44, 10).Value = 36.61125731468201
